$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F ("Trening" = which part of the training a GPS sample
#     belongs to) ------------------------------------------------------
$ws.Cells.Item(1,6).Value = "Trening"

# Give the new header cell the same bold/bordered/centered look as the
# other header cells (copy format only from A1, which already carries it).
$ws.Cells.Item(1,1).Copy()
$ws.Cells.Item(1,6).PasteSpecial(-4122)

# --- Replace / extend the data rows ----------------------------------
# The raw export now has 12 samples (rows 2-13) instead of 6, split
# between "Duża Gra" and "Mała Gra" parts of the session, and the
# Timestamp column is now a real Excel date/time value instead of text.

$ws.Cells.Item(2,1).Value = 45684.59241365741
$ws.Cells.Item(2,2).Value = 584.5
$ws.Cells.Item(2,3).Value = 12.4
$ws.Cells.Item(2,4).Value = 1.947457892554148
$ws.Cells.Item(2,5).Value = "10-15"
$ws.Cells.Item(2,6).Value = "Duża Gra"

$ws.Cells.Item(3,1).Value = 45684.59327361111
$ws.Cells.Item(3,2).Value = 658.8
$ws.Cells.Item(3,3).Value = 14.41
$ws.Cells.Item(3,4).Value = 1.981797490801133
$ws.Cells.Item(3,5).Value = "10-15"
$ws.Cells.Item(3,6).Value = "Duża Gra"

$ws.Cells.Item(4,1).Value = 45684.59395532408
$ws.Cells.Item(4,2).Value = 717.7
$ws.Cells.Item(4,3).Value = 13.39
$ws.Cells.Item(4,4).Value = 2.001448018210273
$ws.Cells.Item(4,5).Value = "10-15"
$ws.Cells.Item(4,6).Value = "Duża Gra"

$ws.Cells.Item(5,1).Value = 45684.59232800926
$ws.Cells.Item(5,2).Value = 577.1
$ws.Cells.Item(5,3).Value = 9.630000000000001
$ws.Cells.Item(5,4).Value = 1.786505733217511
$ws.Cells.Item(5,5).Value = "5-10"
$ws.Cells.Item(5,6).Value = "Duża Gra"

$ws.Cells.Item(6,1).Value = 45684.59309305556
$ws.Cells.Item(6,2).Value = 643.2
$ws.Cells.Item(6,3).Value = 9.630000000000001
$ws.Cells.Item(6,4).Value = 1.814349157469613
$ws.Cells.Item(6,5).Value = "5-10"
$ws.Cells.Item(6,6).Value = "Duża Gra"

$ws.Cells.Item(7,1).Value = 45684.59375625
$ws.Cells.Item(7,2).Value = 700.5
$ws.Cells.Item(7,3).Value = 9.51
$ws.Cells.Item(7,4).Value = 1.814271739551
$ws.Cells.Item(7,5).Value = "5-10"
$ws.Cells.Item(7,6).Value = "Duża Gra"

$ws.Cells.Item(8,1).Value = 45684.59845763889
$ws.Cells.Item(8,2).Value = 1106.7
$ws.Cells.Item(8,3).Value = 13.98
$ws.Cells.Item(8,4).Value = 3.457617555345807
$ws.Cells.Item(8,5).Value = "10-15"
$ws.Cells.Item(8,6).Value = "Mała Gra"

$ws.Cells.Item(9,1).Value = 45684.60006527778
$ws.Cells.Item(9,2).Value = 1245.6
$ws.Cells.Item(9,3).Value = 11.71
$ws.Cells.Item(9,4).Value = 3.00513754572187
$ws.Cells.Item(9,5).Value = "10-15"
$ws.Cells.Item(9,6).Value = "Mała Gra"

$ws.Cells.Item(10,1).Value = 45684.60120532408
$ws.Cells.Item(10,2).Value = 1344.1
$ws.Cells.Item(10,3).Value = 12.98
$ws.Cells.Item(10,4).Value = 3.621854237147737
$ws.Cells.Item(10,5).Value = "10-15"
$ws.Cells.Item(10,6).Value = "Mała Gra"

$ws.Cells.Item(11,1).Value = 45684.5970386574
$ws.Cells.Item(11,2).Value = 984.1
$ws.Cells.Item(11,3).Value = 9.880000000000001
$ws.Cells.Item(11,4).Value = 2.905950750623432
$ws.Cells.Item(11,5).Value = "5-10"
$ws.Cells.Item(11,6).Value = "Mała Gra"

$ws.Cells.Item(12,1).Value = 45684.60120185185
$ws.Cells.Item(12,2).Value = 1343.8
$ws.Cells.Item(12,3).Value = 8.550000000000001
$ws.Cells.Item(12,4).Value = 2.716314515897204
$ws.Cells.Item(12,5).Value = "5-10"
$ws.Cells.Item(12,6).Value = "Mała Gra"

$ws.Cells.Item(13,1).Value = 45684.60276898148
$ws.Cells.Item(13,2).Value = 1479.2
$ws.Cells.Item(13,3).Value = 9.5
$ws.Cells.Item(13,4).Value = 2.758128285408019
$ws.Cells.Item(13,5).Value = "5-10"
$ws.Cells.Item(13,6).Value = "Mała Gra"

# --- Number format for column A (date/time) ---------------------------
# Mirror the author's own workflow: first try a lowercase format code
# (creates numFmt 164, left unused on any cell) then settle on the
# uppercase one (numFmt 165) that is actually applied to every row.
$ws.Cells.Item(2,1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
